$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.264.50"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "1.564.52"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").Value = "'210.70"
$ws.Range("E5").Value = "  +1.14%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "'22.18"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "'0.0872"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").Value = "1.787.36"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.576.25"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "'3.76"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "27.236.86"
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'218.01"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0704"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.45"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'151.60"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("D26").Value = "'6.63"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("D28").Value = "'15.01"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("D31").Value = "'0.0471"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").Value = "1.453.94"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  +4.95%  "
$ws.Range("D36").Value = "'1.63"
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").Value = "'0.543"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D45").Value = "'64.47"
$ws.Range("E45").Value = "  -0.58%  "
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "1.700.19"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "'85.86"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  +1.60%  "
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("D51").Value = "'0.0944"
$ws.Range("E51").Value = "  -1.67%  "
